# Bump the proposal number in the title run from
# "Angebot Nr. 2026-02-12-8" to "Angebot Nr. 2026-02-12-10".
#
# (The rest of the reference diff only swaps the internal r:id
# attributes used by the "Referenzen:" hyperlinks for new
# auto-generated relationship ids while keeping the same hyperlink
# targets/text - a non-semantic artifact of the authoring tool's
# save process, not a content change reachable through the Word
# object model, so it is intentionally left alone here.)

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Angebot Nr. 2026-02-12-8",  # FindText
    $true,                       # MatchCase
    $false,                      # MatchWholeWord
    $false,                      # MatchWildcards
    $false,                      # MatchSoundsLike
    $false,                      # MatchAllWordForms
    $true,                       # Forward
    1,                           # Wrap (wdFindContinue)
    $false,                      # Format
    "Angebot Nr. 2026-02-12-10", # ReplaceWith
    2                            # Replace (wdReplaceAll)
)
